# ETK Key Task README update ("Update KT README for example images"):
#   1. Rename the sample-stimuli image file names referenced in the file
#      list: KeyTaskExample00.png -> KTimage00.png and
#      KeyTaskExample04.png -> KTimage06.png.
#   2. Change "A set of five example stimuli" to "A set of trivial example
#      stimuli" in the intro sentence before the Survey Flow section.
#   3. The Word "_GoBack" bookmark (left over from the last editing
#      position) moves from the end of the "...keyValues." paragraph to
#      sit right after "A set of trivial " in the following paragraph -
#      mirror that by deleting it and re-adding it at the new spot.

$d = $word.ActiveDocument

# --- 1) Sample stimuli image file names -------------------------------
$f = $d.Content
$f.Find.Execute("KeyTaskExample00") | Out-Null
if ($f.Find.Found) {
    $d.Range($f.Start, $f.End).Text = "KTimage00"
}

$f = $d.Content
$f.Find.Execute("KeyTaskExample04") | Out-Null
if ($f.Find.Found) {
    $d.Range($f.Start, $f.End).Text = "KTimage06"
}

# --- 2) "five" -> "trivial" --------------------------------------------
$f = $d.Content
$f.Find.Execute("A set of five example") | Out-Null
if ($f.Find.Found) {
    $d.Range($f.Start, $f.End).Text = "A set of trivial example"
}

# --- 3) Move the _GoBack bookmark --------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$f = $d.Content
$f.Find.Execute("example stimuli are included") | Out-Null
if ($f.Find.Found) {
    $newBookmarkRange = $d.Range($f.Start, $f.Start)
    $d.Bookmarks.Add("_GoBack", $newBookmarkRange) | Out-Null
}

Write-Output "done"
